$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 83
$ws1.Range("F3").Value = 21804
$ws1.Range("F4").Value = 1145
$ws1.Range("F5").Value = 8141
$ws1.Range("F6").Value = 567
$ws1.Range("F7").Value = 329
$ws1.Range("F9").Value = 208
$ws1.Range("F10").Value = 214
$ws1.Range("F11").Value = 548
$ws1.Range("F13").Value = 249
$ws1.Range("F14").Value = 39
$ws1.Range("F15").Value = 1376
$ws1.Range("F16").Value = 572
$ws1.Range("F17").Value = 91
$ws1.Range("F18").Value = 721
$ws1.Range("F19").Value = 61
$ws1.Range("F20").Value = 107
$ws1.Range("F21").Value = 99
$ws1.Range("F22").Value = 377
$ws1.Range("F23").Value = 1241
$ws1.Range("F24").Value = 84
$ws1.Range("F25").Value = 56
$ws1.Range("F26").Value = 246
$ws1.Range("F28").Value = 624
$ws1.Range("F30").Value = 171
$ws1.Range("F31").Value = 5276
$ws1.Range("F32").Value = 39
$ws1.Range("F34").Value = 76
$ws1.Range("F35").Value = 62
$ws1.Range("F36").Value = 13536
$ws1.Range("F37").Value = 1388
$ws1.Range("F38").Value = 164
$ws1.Range("F39").Value = 65
$ws1.Range("F41").Value = 354
$ws1.Range("F42").Value = 485
$ws1.Range("F43").Value = 4105
$ws1.Range("F44").Value = 53

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 83
$ws4.Range("F3").Value = 21804
$ws4.Range("F4").Value = 1145
$ws4.Range("F5").Value = 8141
$ws4.Range("F6").Value = 567
$ws4.Range("F7").Value = 329
$ws4.Range("F9").Value = 208
$ws4.Range("F10").Value = 214
$ws4.Range("F11").Value = 548
$ws4.Range("F13").Value = 249
$ws4.Range("F14").Value = 41
$ws4.Range("F15").Value = 1376
$ws4.Range("F16").Value = 572
$ws4.Range("F17").Value = 91
$ws4.Range("F18").Value = 721
$ws4.Range("F19").Value = 61
$ws4.Range("F20").Value = 107
$ws4.Range("F21").Value = 99
$ws4.Range("F22").Value = 377
$ws4.Range("F23").Value = 1241
$ws4.Range("F24").Value = 84
$ws4.Range("F25").Value = 56
$ws4.Range("F26").Value = 246
$ws4.Range("F29").Value = 624
$ws4.Range("F32").Value = 171
$ws4.Range("F34").Value = 5276
$ws4.Range("F35").Value = 39
$ws4.Range("F37").Value = 76
$ws4.Range("F38").Value = 62
$ws4.Range("F39").Value = 13536
$ws4.Range("F40").Value = 1388
$ws4.Range("F41").Value = 164
$ws4.Range("F42").Value = 65
$ws4.Range("F44").Value = 354
$ws4.Range("F45").Value = 486
$ws4.Range("F46").Value = 4105
$ws4.Range("F47").Value = 53
